# Sprint 7.1 Observaciones - "se agregar codigo de barras" commit
#
# Summary of the change being applied:
#  - A new worksheet "Hoja2" is added (after "Hoja1") holding two rows of
#    barcode/scale lookup data plus a couple of totals formulas.
#  - On "Hoja1", the cash-register ("ticket de cierre de cajas") block gets
#    real figures (column C) and two running-total formulas, and the labels
#    for "Monto Total de Retiros/Cancelaciones" vs "Saldo Total/Saldo
#    Caja:Fisico" are re-ordered.
#  - The typo "tarje" -> "tarjeta" is fixed in the observations text.
#  - A few highlight colors (green/yellow) are touched up, column B is
#    widened, and the active selection moves to B22.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Add the "Hoja2" worksheet right after "Hoja1"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Hoja2"

# Row 13 - "menudeo" style scale lookup
$ws2.Range("B13").Value = 2
$ws2.Range("C13").Value = 2
$ws2.Range("D13").Value = 478
$ws2.Range("E13").Value = 3
$ws2.Range("F13").Value = 0
$ws2.Range("G13").Value = 1650
$ws2.Range("H13").Value = 3
$ws2.Range("I13").Value = 6
$ws2.Range("J13").Value = 550
$ws2.Range("K13").Value = 467
$ws2.Range("L13").Value = 0
$ws2.Range("M13").Value = 550
$ws2.Range("N13").Value = 0
$ws2.Range("O13").Value = "NULL"
$ws2.Range("P13").Value = 0
$ws2.Range("Q13").Value = 41.25

# Row 15 - "mayoreo" style scale lookup
$ws2.Range("B15").Value = 2
$ws2.Range("C15").Value = 2
$ws2.Range("D15").Value = 478
$ws2.Range("E15").Value = 3
$ws2.Range("F15").Value = 0
$ws2.Range("G15").Value = 1000
$ws2.Range("H15").Value = 3
$ws2.Range("I15").Value = 6
$ws2.Range("J15").Value = 550
$ws2.Range("K15").Value = 467
$ws2.Range("L15").Value = 0
$ws2.Range("M15").Value = 550
$ws2.Range("N15").Value = 0
$ws2.Range("O15").Value = "NULL"
$ws2.Range("P15").Value = 0
$ws2.Range("Q15").Formula = "=G15*0.025"

# Row 16 - commission subtotal
$ws2.Range("Q16").Formula = "=Q13+Q15"

# Row 19 - grand totals
$ws2.Range("G19").Formula = "=G13+G15"
$ws2.Range("H19").Formula = "=G19*0.025"

$ws2.Range("Q17").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Fix the "tarje" -> "tarjeta" typo on Hoja1!B23 (must run AFTER the
#    "NULL" strings above so the shared-string table lands on the same
#    indices as the target file).
# ---------------------------------------------------------------------
$ws1.Range("B23").Value = "generar ticket si es cobro con tarjeta se regresa todo monto + comision "

# ---------------------------------------------------------------------
# 3) Re-label / reorder the totals block (rows 27-30) and fill in the
#    new figures + formulas in column C (rows 24-30).
# ---------------------------------------------------------------------
$ws1.Range("B27").Value = "Monto Total de Retiros"
$ws1.Range("B28").Value = "Monto Total de Cancelaciones"
$ws1.Range("B29").Value = "Saldo Total"
$ws1.Range("B30").Value = "Saldo Caja:Fisico"

$ws1.Range("C24").Value = 20
$ws1.Range("C25").Value = 1500
$ws1.Range("C26").Value = 1000
$ws1.Range("C27").Value = 300
$ws1.Range("C28").Value = 100
$ws1.Range("C29").Formula = "=C25+C26-C28"
$ws1.Range("C30").Formula = "=C29-C26-C27"

# Give the new column-C block a plain "no fill" style override (same family
# as the rest of the block, just explicitly applied) so it reads as its own
# style, matching the new cellXfs entry introduced for this block.
$ws1.Range("C24:C30").Interior.ColorIndex = -4142

# ---------------------------------------------------------------------
# 4) Highlight colors. Grab "clean" donor cells BEFORE anything recolors
#    them, so every paste uses the correct source color.
# ---------------------------------------------------------------------
$yellowSrc = $ws1.Range("Q7")   # untouched yellow (fillId 3) cell
$greenSrc  = $ws1.Range("Q2")   # untouched green  (fillId 2 / theme) cell

# B22 -> yellow
$yellowSrc.Copy() | Out-Null
$ws1.Range("B22").PasteSpecial(-4122) | Out-Null

# B12 -> yellow
$yellowSrc.Copy() | Out-Null
$ws1.Range("B12").PasteSpecial(-4122) | Out-Null

# B23, B24 -> green
$greenSrc.Copy() | Out-Null
$ws1.Range("B23:B24").PasteSpecial(-4122) | Out-Null

# Q19 (new, empty) -> green
$greenSrc.Copy() | Out-Null
$ws1.Range("Q19").PasteSpecial(-4122) | Out-Null

# Q6 switches from yellow to green
$greenSrc.Copy() | Out-Null
$ws1.Range("Q6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5) Column B gets a bit wider and the selection moves to B22.
# ---------------------------------------------------------------------
$ws1.Columns("B").ColumnWidth = 30.140625

$ws1.Activate()
$ws1.Range("B22").Select() | Out-Null
